$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "basket" row (B1 / N-5, qty 1) has been dropped from the BOM, so
# remove row 23 entirely (rows below it shift up: old 24->23, 25->24, 26->25).
# Update the time-stamp note first while it's still on row 25 (it moves to
# row 24 once row 23 is deleted). A leading apostrophe keeps it stored as
# text (matching the cell's existing quote-prefixed style) instead of being
# reinterpreted as a time value.
$ws.Range("E25").Value = "'17:42"

$ws.Rows(23).Delete()
